$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 123; this shifts existing rows 123-129 down to 124-130
$ws.Rows.Item(123).Insert()

# Populate the new row 123 with the new weekly price record
$ws.Range("A123").Value = 10
$ws.Range("B123").Value = "Vega Modelo de Temuco"
$ws.Range("C123").Value = "La Araucanía"
$ws.Range("D123").Value = 44826
$ws.Range("E123").Value = 9
$ws.Range("F123").Value = 100114002
$ws.Range("G123").Value = "Camote"
$ws.Range("H123").Value = "Sin especificar"
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 100
$ws.Range("K123").Value = 20000
$ws.Range("L123").Value = 20000
$ws.Range("M123").Value = 20000
$ws.Range("N123").Value = "$/malla 20 kilos"
$ws.Range("O123").Value = "Perú"
$ws.Range("P123").Value = 1000
$ws.Range("Q123").Value = 20
$ws.Range("R123").Value = "Hortaliza"
